$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RejectedHistoryExport")

$ws.Range("A2").Value = "Pradip A. Dey"
$ws.Range("C2").Value = "Grade 3"
$ws.Range("E2").Value = "26-03-2025"
$ws.Range("H2").Value = "Student requires early pickup due to health concerns"

$ws.Columns.Item(8).ColumnWidth = 46.3
